$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 values
$ws.Range("A4").Value = -0.1119805975156518
$ws.Range("B4").Value = 1.3882526130365052
$ws.Range("C4").Value = -0.45866025557780671
$ws.Range("D4").Value = 1.538582730249298
$ws.Range("E4").Value = -0.38809713933506723
$ws.Range("F4").Value = -1.4572817484913592
$ws.Range("G4").Value = 1.4572817484913592

# Update row 5 values
$ws.Range("A5").Value = -0.51388356394168977
$ws.Range("B5").Value = 1.5462526341887264
$ws.Range("C5").Value = -0.89124283776155755
$ws.Range("D5").Value = 1.1535535524900022
$ws.Range("E5").Value = -1.3038836697027949
$ws.Range("F5").Value = -1.3161555160058802
$ws.Range("G5").Value = 1.0231651855197226

# Update sheet view: scroll back to top-left A1 (remove topLeftCell="A21"), and select A4:G4
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A4:G4").Select()
